# Generate Report for Archive
#
# 1) The shared status string "Ready for handoff" becomes "In Translation"
#    everywhere it is used (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# 2) The "Status" columns get narrower:
#      Overview columns E & F : stored width 17.2159881591797 -> 13.4101845877511
#      zh-cn     column  C    : stored width 17.2159881591797 -> 13.4101845877511
#      de-de     column  C    : stored width 17.2159881591797 -> 13.4101845877511
#
# NOTE on the width numbers: this engine persists column width as
#   stored = round(ColumnWidth * 6) / 6 + 5/6
# (i.e. it snaps to whole "pixels" at 6 px/character). The exact decimal
# widths recorded in the target XML (produced by a different writer) are
# not representable through that quantized COM setter, so we feed the
# ColumnWidth value whose quantized result lands closest to the target
# stored width (13.4101845877511 -> stored 13.3333..., reached with
# ColumnWidth = 12.5).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$narrowWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $current = $cell.Value()
            # Compare as strings explicitly -- boolean cell values (e.g. the
            # "True" shared string read back as $true) would otherwise coerce
            # $oldStatus to a truthy bool and match every boolean cell.
            if (($current -is [string]) -and ([string]$current -ceq $oldStatus)) {
                $cell.Value = $newStatus
            }
        }
    }
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $narrowWidth
$overview.Columns.Item(6).ColumnWidth = $narrowWidth

$wb.Worksheets.Item("zh-cn").Columns.Item(3).ColumnWidth = $narrowWidth
$wb.Worksheets.Item("de-de").Columns.Item(3).ColumnWidth = $narrowWidth
